$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.855.08"
$ws.Range("E2").Value = "  +2.93%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.562.67"
$ws.Range("E3").Value = "  +2.08%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "614.28"
$ws.Range("E5").Value = "  +6.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.48"
$ws.Range("E6").Value = "  +1.17%  "
$ws.Range("E7").Value = "  +2.32%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.558.40"
$ws.Range("E8").Value = "  +2.14%  "
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("E10").Value = "  +5.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.07"
$ws.Range("E11").Value = "  +8.62%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.586"
$ws.Range("E12").Value = "  +1.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "46.87"
$ws.Range("E13").Value = "  +1.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000279"
$ws.Range("E14").Value = "  +2.90%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.132.67"
$ws.Range("E15").Value = "  +2.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.41"
$ws.Range("E16").Value = "  -0.68%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "620.20"
$ws.Range("E17").Value = "  -0.61%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.968.35"
$ws.Range("E18").Value = "  +3.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.559.60"
$ws.Range("E19").Value = "  +1.91%  "
$ws.Range("E20").Value = "  -0.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.40"
$ws.Range("E21").Value = "  +1.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.884"
$ws.Range("E22").Value = "  +0.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.52"
$ws.Range("E23").Value = "  -13.76%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.77"
$ws.Range("E24").Value = "  -0.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "96.86"
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.82"
$ws.Range("E26").Value = "  +0.64%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.52"
$ws.Range("E29").Value = "  +3.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.11"
$ws.Range("E30").Value = "  -1.72%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.53"
$ws.Range("E31").Value = "  +0.93%  "
$ws.Range("E32").Value = "  -1.87%  "
$ws.Range("B33").Value = "Mantle"
$ws.Range("C33").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.31"
$ws.Range("E33").Value = "  -0.25%  "
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.02"
$ws.Range("E34").Value = "  +0.54%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "574.13"
$ws.Range("E35").Value = "  -7.36%  "
$ws.Range("E36").Value = "  +5.19%  "
$ws.Range("E37").Value = "  -0.85%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "10.84"
$ws.Range("E38").Value = "  +1.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "57.80"
$ws.Range("E39").Value = "  +2.42%  "
$ws.Range("E40").Value = "  +6.80%  "
$ws.Range("E41").Value = "  +0.20%  "
$ws.Range("E42").Value = "  +5.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.363.59"
$ws.Range("E43").Value = "  +0.41%  "
$ws.Range("B44").Value = "ThetaToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.05"
$ws.Range("E44").Value = "  +9.92%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.322"
$ws.Range("E45").Value = "  -1.25%  "
$ws.Range("B46").Value = "PEPE"
$ws.Range("C46").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₃0711"
$ws.Range("E46").Value = "  +3.17%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "33.15"
$ws.Range("E47").Value = "  +1.60%  "
$ws.Range("E48").Value = "  +3.75%  "
$ws.Range("E49").Value = "  +0.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "134.13"
$ws.Range("E50").Value = "  +2.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.70"
$ws.Range("E51").Value = "  +2.43%  "
